# Auto-generated Excel COM-interop script
# Applies numeric refresh of market/profit data across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 11
$ws.Range("H11").Value = 2988.7844
$ws.Range("I11").Value = 2988.7844
$ws.Range("K11").Value = 2988.7844
$ws.Range("M11").Value = -2848.7844

# Row 15
$ws.Range("H15").Value = 1248.746
$ws.Range("I15").Value = 1248.746
$ws.Range("K15").Value = 3746.238
$ws.Range("M15").Value = -3577.238

# Row 18
$ws.Range("H18").Value = 2329.6667
$ws.Range("I18").Value = 2329.6667
$ws.Range("K18").Value = 2329.6667
$ws.Range("M18").Value = -2045.6667

# Row 46
$ws.Range("H46").Value = 7999
$ws.Range("J46").Value = 7999
$ws.Range("L46").Value = 23997
$ws.Range("N46").Value = -24235

# Row 48
$ws.Range("H48").Value = 867
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

# Row 56
$ws.Range("H56").Value = 867
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

# Row 60
$ws.Range("H60").Value = 7999
$ws.Range("J60").Value = 7999
$ws.Range("L60").Value = 23997
$ws.Range("N60").Value = -24965

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# Row 98
$ws.Range("H98").Value = 995.5789
$ws.Range("I98").Value = 759.94116
$ws.Range("K98").Value = 759.94116
$ws.Range("M98").Value = 738.05884

# Row 112
$ws.Range("H112").Value = 3248149
$ws.Range("I112").Value = 748.5
$ws.Range("J112").Value = 3789382.5
$ws.Range("K112").Value = 2245.5
$ws.Range("L112").Value = 11368147.5
$ws.Range("M112").Value = -1137.5
$ws.Range("N112").Value = -11370363.5

# Row 113
$ws.Range("H113").Value = 55215.105
$ws.Range("I113").Value = 2698.5
$ws.Range("J113").Value = 145243.58
$ws.Range("K113").Value = 2698.5
$ws.Range("L113").Value = 145243.58
$ws.Range("M113").Value = 555.5
$ws.Range("N113").Value = -151751.58

# Row 122
$ws.Range("H122").Value = 995.5789
$ws.Range("I122").Value = 759.94116
$ws.Range("K122").Value = 2279.82348
$ws.Range("M122").Value = 170.17652

# Row 127
$ws.Range("H127").Value = 1321.9445
$ws.Range("I127").Value = 1321.9445
$ws.Range("K127").Value = 3965.8335
$ws.Range("M127").Value = 994.1664999999998

# Row 131
$ws.Range("H131").Value = 3076.2222
$ws.Range("I131").Value = 1612.5294
$ws.Range("J131").Value = 5564.5
$ws.Range("K131").Value = 4837.5882
$ws.Range("L131").Value = 16693.5
$ws.Range("M131").Value = 202.4117999999999
$ws.Range("N131").Value = -26773.5

# Row 135
$ws.Range("H135").Value = 1331.2812
$ws.Range("I135").Value = 779.25
$ws.Range("K135").Value = 7013.25
$ws.Range("M135").Value = -4478.25

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 45612.82
$ws.Range("I32").Value = 47506.152
$ws.Range("K32").Value = 47506.152
$ws.Range("M32").Value = -47219.152

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 2545.3333
$ws.Range("I20").Value = 2396.9092
$ws.Range("K20").Value = 2396.9092
$ws.Range("M20").Value = -2149.9092

# Row 86
$ws.Range("H86").Value = 168832.83
$ws.Range("J86").Value = 252362.38
$ws.Range("L86").Value = 252362.38
$ws.Range("N86").Value = -254608.38

# Row 89
$ws.Range("H89").Value = 168832.83
$ws.Range("J89").Value = 252362.38
$ws.Range("L89").Value = 1261811.9
$ws.Range("N89").Value = -1273043.9

# Row 132
$ws.Range("H132").Value = 110709
$ws.Range("J132").Value = 110709
$ws.Range("L132").Value = 110709
$ws.Range("N132").Value = -120829

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

# Row 50
$ws.Range("H50").Value = 37598.6
$ws.Range("J50").Value = 37598.6
$ws.Range("L50").Value = 37598.6
$ws.Range("N50").Value = -38848.6

# Row 51
$ws.Range("H51").Value = 66843.75
$ws.Range("J51").Value = 74110
$ws.Range("L51").Value = 74110
$ws.Range("N51").Value = -75582

# Row 59
$ws.Range("H59").Value = 51642.57
$ws.Range("J59").Value = 56299.6
$ws.Range("L59").Value = 56299.6
$ws.Range("N59").Value = -58589.6

# Row 60
$ws.Range("H60").Value = 11682.667
$ws.Range("J60").Value = 13150.25
$ws.Range("L60").Value = 13150.25
$ws.Range("N60").Value = -14172.25

# Row 61
$ws.Range("H61").Value = 66843.75
$ws.Range("J61").Value = 74110
$ws.Range("L61").Value = 74110
$ws.Range("N61").Value = -74806

# Row 141
$ws.Range("H141").Value = 327334.6
$ws.Range("J141").Value = 327334.6
$ws.Range("L141").Value = 327334.6
$ws.Range("N141").Value = -337694.6

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 107
$ws.Range("H107").Value = 1624.75
$ws.Range("I107").Value = 616.3333
$ws.Range("J107").Value = 1857.4615
$ws.Range("K107").Value = 1848.9999
$ws.Range("L107").Value = 5572.3845
$ws.Range("M107").Value = 71.00009999999997
$ws.Range("N107").Value = -9412.3845

# Row 116
$ws.Range("H116").Value = 3666
$ws.Range("I116").Value = 998
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 2994
$ws.Range("L116").Value = 15000
$ws.Range("M116").Value = 448
$ws.Range("N116").Value = -21884

# Row 121
$ws.Range("H121").Value = 15153180
$ws.Range("I121").Value = 1272
$ws.Range("K121").Value = 3816
$ws.Range("M121").Value = -2506

# Row 131
$ws.Range("H131").Value = 22228018
$ws.Range("J131").Value = 6908.3335
$ws.Range("L131").Value = 20725.0005
$ws.Range("N131").Value = -30805.0005

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

# Row 80
$ws.Range("H80").Value = 3712.611
$ws.Range("I80").Value = 2966.7273
$ws.Range("J80").Value = 4884.7144
$ws.Range("K80").Value = 2966.7273
$ws.Range("L80").Value = 4884.7144
$ws.Range("M80").Value = -1968.7273
$ws.Range("N80").Value = -6880.7144

# Row 83
$ws.Range("H83").Value = 3712.611
$ws.Range("I83").Value = 2966.7273
$ws.Range("J83").Value = 4884.7144
$ws.Range("K83").Value = 14833.6365
$ws.Range("L83").Value = 24423.572
$ws.Range("M83").Value = -9841.636500000001
$ws.Range("N83").Value = -34407.572

# Row 122
$ws.Range("H122").Value = 4662.7085
$ws.Range("I122").Value = 2085.3333
$ws.Range("J122").Value = 8958.333000000001
$ws.Range("K122").Value = 6255.999899999999
$ws.Range("L122").Value = 26874.999
$ws.Range("M122").Value = -3805.999899999999
$ws.Range("N122").Value = -31774.999

# Row 132
$ws.Range("H132").Value = 6990.4443
$ws.Range("I132").Value = 5423.5835
$ws.Range("J132").Value = 8243.933999999999
$ws.Range("K132").Value = 16270.7505
$ws.Range("L132").Value = 24731.802
$ws.Range("M132").Value = -13740.7505
$ws.Range("N132").Value = -29791.802

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 56
$ws.Range("H56").Value = 14247.25
$ws.Range("I56").Value = 6995.6
$ws.Range("K56").Value = 6995.6
$ws.Range("M56").Value = -6304.6

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 64
$ws.Range("H64").Value = 20000
$ws.Range("I64").Value = 20000
$ws.Range("K64").Value = 20000
$ws.Range("M64").Value = -19752

# Row 67
$ws.Range("H67").Value = 20000
$ws.Range("I67").Value = 20000
$ws.Range("K67").Value = 20000
$ws.Range("M67").Value = -19142

# Row 107
$ws.Range("H107").Value = 1274.2433
$ws.Range("I107").Value = 967.96
$ws.Range("K107").Value = 2903.88
$ws.Range("M107").Value = -983.8800000000001

Write-Output "Applied scheduled market-data refresh to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR."
